$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: title "COP CAR 2" becomes "COP 2" (shared string reused, same index),
# and the year in B1 changes to 1900.
$ws.Range("A1").Value = "COP 2"
$ws.Range("B1").Value = 1900
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = "James"

# Row 2: gains the rest of its record (it previously only had B2).
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = "Hames"

# Old row 3 ("COP CAR" / 1900 / 2 / 3 / James) is removed entirely.
$ws.Range("A3:E3").ClearContents()

# New row 4 is added with a fresh record (set E4 before A4 so the shared
# strings table picks up "Jake" ahead of "COPP", matching save order).
$ws.Range("E4").Value = "Jake"
$ws.Range("A4").Value = "COPP"
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 13

# Leave the saved selection where the author left it.
$ws.Range("E9").Select()
